$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 27; $r -le 66; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
